# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) wherever they appear as "System, <email>",
# turning them into "<email>, System".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    if ($cell.Value -eq $oldText) {
        $cell.Value = $newText
    }
}
